$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.170.33"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +3.32%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.118.08"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.30%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.27%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "544.85"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.56%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.97"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +6.48%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.21%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.107.92"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +3.75%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.500"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +4.19%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +4.15%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +3.05%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.85%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000228"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +8.85%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.13"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.58%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.616.65"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +3.10%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.183.25"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.95%  "

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.49%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.108.42"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.43%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.75"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +4.34%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "488.69"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.22%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +3.35%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.708"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +3.87%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.21"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +4.21%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.68"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +4.77%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.39"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +4.06%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.66%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +3.84%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.23"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.00%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.41%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "26.56"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.49%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.93"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.64%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.16"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +4.33%  "

$ws.Range("B33").Value = "Stacks"
$ws.Range("C33").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.39"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.20%  "

$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "57.59"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.37%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "505.24"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.31%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.40"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +7.81%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.07"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +5.40%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.279.06"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +8.39%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +4.70%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0803"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +4.26%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.120"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +3.77%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.73"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +7.43%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.18"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +4.08%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.259"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +5.38%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.06%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +5.06%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "123.83"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.67%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +11.85%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.03"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +5.91%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +4.51%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.42"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +3.71%  "
